# Apply the "Mon Dec  4 08:14:30 UTC 2023" cryptos-list refresh (GitHub Actions scrape).
#
# Every value in columns D (Price) and E (Volume(1h)) is stored as literal text in
# this workbook (even when it looks like a plain number, e.g. "231.99"), so plain
# `.Value = "231.92"` assignment would let Excel auto-convert it to a real number.
# To keep cells textual (matching the source data / original inlineStr cells) we
# prefix genuinely-numeric-looking Price strings with a leading apostrophe, which
# is the standard Excel "treat as text" input convention.
#
# Rows 28/29 (Monero <-> Kaspa) also swap their Coin name/Link/Price/Volume data;
# the rank index in column A is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '41.469.07'
$ws.Range("E2").Value = '  +5.03%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.245.58'
$ws.Range("E3").Value = '  +3.81%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.16%  '

# Row 5: BNB
$ws.Range("D5").Value = '''231.92'  # was '231.99'; forced text (looks numeric)
$ws.Range("E5").Value = '  +1.74%  '

# Row 6: XRP
$ws.Range("D6").Value = '''0.638'  # was '0.637'; forced text (looks numeric)
$ws.Range("E6").Value = '  +2.29%  '

# Row 7: Solana
$ws.Range("D7").Value = '''64.00'  # was '63.86'; forced text (looks numeric)
$ws.Range("E7").Value = '  -0.12%  '

# Row 8: USDC
$ws.Range("E8").Value = '  +0.11%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  +3.02%  '

# Row 10: OKB
$ws.Range("D10").Value = '''59.65'  # was '59.41'; forced text (looks numeric)
$ws.Range("E10").Value = '  +2.63%  '

# Row 11: Dogecoin
$ws.Range("D11").Value = '''0.0902'  # was '0.0903'; forced text (looks numeric)
$ws.Range("E11").Value = '  +5.42%  '

# Row 12: TRON
$ws.Range("E12").Value = '  +0.95%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '2.582.68'
$ws.Range("E13").Value = '  +3.96%  '

# Row 14: Chainlink
$ws.Range("D14").Value = '''16.13'  # was '16.12'; forced text (looks numeric)
$ws.Range("E14").Value = '  -0.83%  '

# Row 15: Avalanche
$ws.Range("D15").Value = '''22.56'  # was '22.54'; forced text (looks numeric)
$ws.Range("E15").Value = '  +1.76%  '

# Row 16: Polygon
$ws.Range("D16").Value = '''0.824'  # was '0.825'; forced text (looks numeric)
$ws.Range("E16").Value = '  +1.16%  '

# Row 17: Polkadot
$ws.Range("D17").Value = '''5.66'  # was '5.65'; forced text (looks numeric)
$ws.Range("E17").Value = '  +2.22%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '2.249.34'
$ws.Range("E18").Value = '  +4.19%  '

# Row 19: WrappedBTC
$ws.Range("D19").Value = '41.354.58'
$ws.Range("E19").Value = '  +4.73%  '

# Row 20: ShibaInu
$ws.Range("D20").Value = '0.0₃0944'
$ws.Range("E20").Value = '  +10.70%  '

# Row 21: Litecoin
$ws.Range("D21").Value = '''73.62'  # was '73.56'; forced text (looks numeric)
$ws.Range("E21").Value = '  +2.29%  '

# Row 22: Uniswap
$ws.Range("D22").Value = '''6.16'  # was '6.18'; forced text (looks numeric)
$ws.Range("E22").Value = '  +0.52%  '

# Row 23: BitcoinCash
$ws.Range("D23").Value = '''252.22'  # was '251.25'; forced text (looks numeric)
$ws.Range("E23").Value = '  +9.68%  '

# Row 24: Dai
$ws.Range("E24").Value = '  -0.04%  '

# Row 25: PancakeSwap
$ws.Range("E25").Value = '  +1.73%  '

# Row 26: Toncoin
$ws.Range("D26").Value = '''2.33'  # was '2.32'; forced text (looks numeric)
$ws.Range("E26").Value = '  +1.18%  '

# Row 27: Cosmos
$ws.Range("D27").Value = '''9.92'  # was '9.90'; forced text (looks numeric)
$ws.Range("E27").Value = '  +2.21%  '

# Row 28: Monero/Kaspa (rows swap)
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.147'  # was '173.20'; forced text (looks numeric)
$ws.Range("E28").Value = '  +3.88%  '

# Row 29: Kaspa/Monero (rows swap)
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '''173.21'  # was '0.145'; forced text (looks numeric)
$ws.Range("E29").Value = '  +0.49%  '

# Row 30: EthereumClassic
$ws.Range("E30").Value = '  +2.61%  '

# Row 31: ImmutableX
$ws.Range("D31").Value = '''1.44'  # was '1.45'; forced text (looks numeric)
$ws.Range("E31").Value = '  +0.22%  '

# Row 32: WEMIXToken
$ws.Range("E32").Value = '  +8.19%  '

# Row 33: Stellar
$ws.Range("E33").Value = '  +1.92%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range("D34").Value = '''5.04'  # was '5.06'; forced text (looks numeric)
$ws.Range("E34").Value = '  +6.23%  '

# Row 35: Filecoin
$ws.Range("D35").Value = '''4.76'  # was '4.74'; forced text (looks numeric)
$ws.Range("E35").Value = '  +3.25%  '

# Row 36: Hedera
$ws.Range("D36").Value = '''0.0638'  # was '0.0634'; forced text (looks numeric)
$ws.Range("E36").Value = '  +3.44%  '

# Row 37: THORChain
$ws.Range("D37").Value = '''6.90'  # was '6.93'; forced text (looks numeric)
$ws.Range("E37").Value = '  -3.23%  '

# Row 38: RenderToken
$ws.Range("D38").Value = '''3.86'  # was '3.83'; forced text (looks numeric)
$ws.Range("E38").Value = '  +8.05%  '

# Row 39: LidoDAOToken
$ws.Range("D39").Value = '''2.44'  # was '2.45'; forced text (looks numeric)
$ws.Range("E39").Value = '  -0.65%  '

# Row 40: TerraClassic
$ws.Range("D40").Value = '''0.000267'  # was '0.000271'; forced text (looks numeric)
$ws.Range("E40").Value = '  +69.91%  '

# Row 41: BinanceUSD
$ws.Range("E41").Value = '  +0.17%  '

# Row 42: FTXToken
$ws.Range("D42").Value = '''4.92'  # was '4.91'; forced text (looks numeric)
$ws.Range("E42").Value = '  +14.88%  '

# Row 43: VeChain
$ws.Range("D43").Value = '''0.0242'  # was '0.0240'; forced text (looks numeric)
$ws.Range("E43").Value = '  +5.10%  '

# Row 44: FraxShare
$ws.Range("D44").Value = '''8.97'  # was '8.86'; forced text (looks numeric)
$ws.Range("E44").Value = '  +13.48%  '

# Row 45: Aave
$ws.Range("D45").Value = '''102.59'  # was '102.76'; forced text (looks numeric)
$ws.Range("E45").Value = '  -0.79%  '

# Row 46: InjectiveProtocol
$ws.Range("D46").Value = '''17.73'  # was '17.65'; forced text (looks numeric)
$ws.Range("E46").Value = '  -0.12%  '

# Row 47: TrustWalletToken
$ws.Range("D47").Value = '''1.24'  # was '1.23'; forced text (looks numeric)
$ws.Range("E47").Value = '  +3.76%  '

# Row 48: Maker
$ws.Range("D48").Value = '1.510.90'
$ws.Range("E48").Value = '  -1.17%  '

# Row 49: Cronos
$ws.Range("D49").Value = '''0.0969'  # was '0.0945'; forced text (looks numeric)
$ws.Range("E49").Value = '  +3.54%  '

# Row 50: ARBITRUM
$ws.Range("E50").Value = '  +1.87%  '

# Row 51: HuobiToken
$ws.Range("D51").Value = '''2.80'  # was '2.79'; forced text (looks numeric)
$ws.Range("E51").Value = '  -0.80%  '
